# Added a concrete curve class and line class
# Tested hermite spline and created a common concrete class for all
# curves to derive from.
#
# This populates Sheet3 with a small trig table (H6:K10), a polyline
# arc-length / parametrisation table (D/E/F columns against the existing
# A/B quarter-circle samples), a couple of stray one-off checks (E18),
# and three new trailer rows (23-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Trig table: angles 0..4*(pi/8) and their 2*cos/2*sin ---------------
$ws.Range("H6").Formula = '=PI()/8'
$ws.Range("I6").Value = 0
$ws.Range("J6").Formula = '=2*COS(I6*$H$6)'
$ws.Range("K6").Formula = '=2*SIN(I6*$H$6)'

$ws.Range("I7").Formula = '=I6+1'
$ws.Range("J7").Formula = '=2*COS(I7*$H$6)'
$ws.Range("K7").Formula = '=2*SIN(I7*$H$6)'

$ws.Range("I8").Formula = '=I7+1'
$ws.Range("J8").Formula = '=2*COS(I8*$H$6)'
$ws.Range("K8").Formula = '=2*SIN(I8*$H$6)'

$ws.Range("I9").Formula = '=I8+1'
$ws.Range("J9").Formula = '=2*COS(I9*$H$6)'
$ws.Range("K9").Formula = '=2*SIN(I9*$H$6)'

$ws.Range("I10").Formula = '=I9+1'
$ws.Range("J10").Formula = '=2*COS(I10*$H$6)'
$ws.Range("K10").Formula = '=2*SIN(I10*$H$6)'

# --- F11 seed value for the arc-length accumulation below ---------------
$ws.Range("F11").Value = 0

# --- D/E/F: chord length, cumulative arc length, and fraction of PI -----
$ws.Range("D12").Formula = '=SQRT((A12-A11)^2+(B12-B11)^2)'
$ws.Range("E12").Formula = '=D12'
$ws.Range("F12").Formula = '=E12/PI()'

$ws.Range("D13").Formula = '=SQRT((A13-A12)^2+(B13-B12)^2)'
$ws.Range("E13").Formula = '=E12+D13'
$ws.Range("F13").Formula = '=E13/PI()'

$ws.Range("D14").Formula = '=SQRT((A14-A13)^2+(B14-B13)^2)'
$ws.Range("E14").Formula = '=E13+D14'
$ws.Range("F14").Formula = '=E14/PI()'

$ws.Range("D15").Formula = '=SQRT((A15-A14)^2+(B15-B14)^2)'
$ws.Range("E15").Formula = '=E14+D15'
$ws.Range("F15").Formula = '=E15/PI()'

$ws.Range("D16").Formula = '=SQRT((A16-A15)^2+(B16-B15)^2)'
$ws.Range("E16").Formula = '=E15+D16'
$ws.Range("F16").Formula = '=E16/PI()'

# --- stray one-off check near the existing 2*SIN(PI()/4) sample ---------
$ws.Range("E18").Formula = '=(1.444-A18)/A18'

# --- trailer rows ---------------------------------------------------------
$ws.Range("B23").Formula = '=SQRT(2)'
$ws.Range("E24").Formula = '=0.25-0.128'
$ws.Range("E25").Formula = '=E24*3.1415'

# --- selection moves to M9 on Sheet3 (matches author's final cursor pos) -
$ws.Range("M9").Select()
